# Update "想去人数" (F column) values across sheets as per upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 782
$ws.Range("F12").Value = 1422
$ws.Range("F15").Value = 1589
$ws.Range("F22").Value = 1508
$ws.Range("F24").Value = 605
$ws.Range("F25").Value = 485
$ws.Range("F31").Value = 2406
$ws.Range("F33").Value = 1346
$ws.Range("F36").Value = 3930

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F29").Value = 26

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 989

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 989
$ws.Range("F8").Value  = 782
$ws.Range("F21").Value = 1422
$ws.Range("F24").Value = 1589
$ws.Range("F30").Value = 1508
$ws.Range("F33").Value = 605
$ws.Range("F34").Value = 485
$ws.Range("F43").Value = 2406
$ws.Range("F48").Value = 1346
$ws.Range("F50").Value = 3930
